# Whitelist update: append 3 new rows (121-123) to the "Whitelist" sheet,
# matching the rows added by the commit "Updated whitelist with broadwayworld.com"
# (secretsanfrancisco.com / sfcv.org / shazam.com follow-on additions).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row data: domain, category, name, city (city left blank/empty, same as row 120)
$newRows = @(
    @("secretsanfrancisco.com", "custom_search", "Brava Theater Center", ""),
    @("sfcv.org",               "music",         "Bing Concert Hall",   ""),
    @("shazam.com",             "music",         "Bing Concert Hall",   "")
)

$startRow = 121
for ($i = 0; $i -lt $newRows.Length; $i++) {
    $r = $startRow + $i
    $rowData = $newRows[$i]

    $ws.Cells.Item($r, 1).Value = $rowData[0]
    $ws.Cells.Item($r, 2).Value = $rowData[1]
    $ws.Cells.Item($r, 3).Value = $rowData[2]

    # The "city" column is an empty string (not a truly blank cell) in the
    # source data, matching every other empty city cell in the sheet
    # (e.g. D120). A plain Value = "" assignment collapses to a blank cell,
    # so force a text entry (leading apostrophe = explicit text) which
    # preserves the empty string, then strip the quote-prefix style it adds
    # so the cell formatting matches its neighbours.
    $ws.Cells.Item($r, 4).Value = "'"
    $ws.Cells.Item($r, 4).Style = "Normal"
}

# Keep the "numbers stored as text" ignored-error annotation in sync with the
# sheet's new used range (best effort - matches how Excel re-flows
# ignoredErrors sqref when the whitelisted text range grows).
$lastRow = $startRow + $newRows.Length - 1
$ws.Range("A1:D$lastRow").Errors.Item(9).Ignore = $true

Write-Host "Added rows 121-123 to Whitelist sheet"
